# Generate Report for Handback
# Update the handoff/handback timestamps recorded on the per-language
# status sheets ("zh-cn" and "de-de") to reflect the freshly generated
# report run.

$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("D2").Value = "2016-02-22 03:26:58"
$wsZh.Range("G2").Value = "2016-02-22 03:27:40"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("D2").Value = "2016-02-22 03:27:10"
$wsDe.Range("G2").Value = "2016-02-22 03:28:00"
